# sn: update itas forms
# Bump the ITAS FTS result form from V3 to V4:
#  - survey sheet: rename the "begin repeat" group from sn_lf_f_2405_v3 to sn_lf_f_2405_v4
#  - survey sheet: add a QR-code format constraint + message to the d_manual_code_id question
#  - settings sheet: bump form_title / form_id to the V4 values

$wb = $excel.ActiveWorkbook

$survey = $wb.Worksheets.Item("survey")
$settings = $wb.Worksheets.Item("settings")

# --- survey!B9 : begin-repeat group name, V3 -> V4 ---
$survey.Range("B9").Value = "sn_lf_f_2405_v4"

# --- survey!F14:G14 : new constraint + constraint_message for d_manual_code_id ---
$survey.Range("F14").Value = 'if(${d_eu_name} = ''SARAYA'', regex(., ''^(SENITAS)\d{4}$''), true())'
$survey.Range("G14").Value = "Le format du QR Code est incorrect. Exemple SENITAS1234"
$survey.Rows.Item(14).RowHeight = 51

# --- settings!A2:B2 : form_title / form_id, V3 -> V4 ---
$settings.Range("A2").Value = "(2024 Mai) 2. ITAS - Formulaire Résultat FTS V4"
$settings.Range("B2").Value = "sn_lf_itas_20305_2_fts_v4"

# --- active-sheet / selection bookkeeping matching the author's last click ---
$settings.Activate()
$settings.Range("B2").Select() | Out-Null
